$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1272.0588
$ws.Range("I41").Value = 1751.6364
$ws.Range("J41").Value = 392.83334
$ws.Range("K41").Value = 1751.6364
$ws.Range("L41").Value = 392.83334
$ws.Range("M41").Value = -1311.6364
$ws.Range("N41").Value = -1272.83334

$ws.Range("H70").Value = 114676.89
$ws.Range("I70").Value = 2096
$ws.Range("J70").Value = 146842.86
$ws.Range("K70").Value = 6288
$ws.Range("L70").Value = 440528.58
$ws.Range("M70").Value = -6018
$ws.Range("N70").Value = -441068.58

$ws.Range("H73").Value = 114676.89
$ws.Range("I73").Value = 2096
$ws.Range("J73").Value = 146842.86
$ws.Range("K73").Value = 6288
$ws.Range("L73").Value = 440528.58
$ws.Range("M73").Value = -5352
$ws.Range("N73").Value = -442400.58

$ws.Range("H132").Value = 1264.4849
$ws.Range("I132").Value = 1210.25
$ws.Range("K132").Value = 3630.75
$ws.Range("M132").Value = -1100.75

$ws.Range("H135").Value = 2398.111
$ws.Range("I135").Value = 1763.5
$ws.Range("J135").Value = 3667.3333
$ws.Range("K135").Value = 15871.5
$ws.Range("L135").Value = 33005.9997
$ws.Range("M135").Value = -13336.5
$ws.Range("N135").Value = -38075.9997

$ws.Range("H138").Value = 6175869
$ws.Range("J138").Value = 6947646.5
$ws.Range("L138").Value = 20842939.5
$ws.Range("N138").Value = -20853219.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 824
$ws.Range("I2").Value = 824
$ws.Range("K2").Value = 824
$ws.Range("M2").Value = -711

$ws.Range("H116").Value = 824
$ws.Range("I116").Value = 824
$ws.Range("K116").Value = 824
$ws.Range("M116").Value = 1470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 824
$ws.Range("I3").Value = 824
$ws.Range("K3").Value = 824
$ws.Range("M3").Value = -710

$ws.Range("H58").Value = 31215
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 31215
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 31215
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -31803

$ws.Range("H75").Value = 22000.334
$ws.Range("I75").Value = 8000.5
$ws.Range("J75").Value = 50000
$ws.Range("K75").Value = 8000.5
$ws.Range("L75").Value = 50000
$ws.Range("M75").Value = -7064.5
$ws.Range("N75").Value = -51872

$ws.Range("H78").Value = 22000.334
$ws.Range("I78").Value = 8000.5
$ws.Range("J78").Value = 50000
$ws.Range("K78").Value = 24001.5
$ws.Range("L78").Value = 150000
$ws.Range("M78").Value = -19321.5
$ws.Range("N78").Value = -159360

$ws.Range("H134").Value = 2009.2134
$ws.Range("J134").Value = 2427.1667
$ws.Range("L134").Value = 7281.500100000001
$ws.Range("N134").Value = -12351.5001

$ws.Range("H138").Value = 69198
$ws.Range("J138").Value = 69198
$ws.Range("L138").Value = 69198
$ws.Range("N138").Value = -79478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 96677.25
$ws.Range("J52").Value = 99000
$ws.Range("L52").Value = 99000
$ws.Range("N52").Value = -99588

$ws.Range("H122").Value = 1333.6
$ws.Range("I122").Value = 1008.3
$ws.Range("J122").Value = 1984.2
$ws.Range("K122").Value = 3024.9
$ws.Range("L122").Value = 5952.6
$ws.Range("M122").Value = -574.8999999999996
$ws.Range("N122").Value = -10852.6

$ws.Range("H132").Value = 3623.2144
$ws.Range("I132").Value = 3444.7144
$ws.Range("J132").Value = 4158.7144
$ws.Range("K132").Value = 10334.1432
$ws.Range("L132").Value = 12476.1432
$ws.Range("M132").Value = -7804.143199999999
$ws.Range("N132").Value = -17536.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 547
$ws.Range("I5").Value = 499.3158
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1497.9474
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1385.9474
$ws.Range("N5").Value = -3224

$ws.Range("H70").Value = 6528.5713
$ws.Range("I70").Value = 6283.3335
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 18850.0005
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = -18535.0005
$ws.Range("N70").Value = -24630

$ws.Range("H73").Value = 6528.5713
$ws.Range("I73").Value = 6283.3335
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 18850.0005
$ws.Range("L73").Value = 24000
$ws.Range("M73").Value = -17758.0005
$ws.Range("N73").Value = -26184

$ws.Range("H75").Value = 1870
$ws.Range("I75").Value = 1560
$ws.Range("J75").Value = 2180
$ws.Range("K75").Value = 4680
$ws.Range("L75").Value = 6540
$ws.Range("M75").Value = -3682
$ws.Range("N75").Value = -8536

$ws.Range("H78").Value = 1870
$ws.Range("I78").Value = 1560
$ws.Range("J78").Value = 2180
$ws.Range("K78").Value = 14040
$ws.Range("L78").Value = 19620
$ws.Range("M78").Value = -9048
$ws.Range("N78").Value = -29604

$ws.Range("H109").Value = 2875.5
$ws.Range("I109").Value = 2450.6
$ws.Range("K109").Value = 7351.799999999999
$ws.Range("M109").Value = -6311.799999999999

$ws.Range("H131").Value = 22225.18
$ws.Range("I131").Value = 101210.2
$ws.Range("J131").Value = 2478.925
$ws.Range("K131").Value = 303630.6
$ws.Range("L131").Value = 7436.775000000001
$ws.Range("M131").Value = -298590.6
$ws.Range("N131").Value = -17516.775

$ws.Range("H134").Value = 6619.875
$ws.Range("I134").Value = 3810.7273
$ws.Range("K134").Value = 11432.1819
$ws.Range("M134").Value = -6362.1819

$ws.Range("H135").Value = 547
$ws.Range("I135").Value = 499.3158
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 4493.8422
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -1958.8422
$ws.Range("N135").Value = -14070

$ws.Range("H137").Value = 3025.2727
$ws.Range("J137").Value = 4241.5
$ws.Range("L137").Value = 12724.5
$ws.Range("N137").Value = -22924.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6571.4287
$ws.Range("I5").Value = 3500
$ws.Range("J5").Value = 7800
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 7800
$ws.Range("M5").Value = -3388
$ws.Range("N5").Value = -8024

$ws.Range("H122").Value = 6050.6665
$ws.Range("I122").Value = 5649.3335
$ws.Range("K122").Value = 16948.0005
$ws.Range("M122").Value = -14498.0005

$ws.Range("H132").Value = 2182.8628
$ws.Range("I132").Value = 1811
$ws.Range("K132").Value = 5433
$ws.Range("M132").Value = -2903

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5684
$ws.Range("I2").Value = 3525
$ws.Range("K2").Value = 3525
$ws.Range("M2").Value = -3413

$ws.Range("H61").Value = 2370.05
$ws.Range("I61").Value = 2118.875
$ws.Range("K61").Value = 2118.875
$ws.Range("M61").Value = -1916.875

$ws.Range("H100").Value = 1800
$ws.Range("J100").Value = 1950
$ws.Range("L100").Value = 1950
$ws.Range("N100").Value = -3032

$ws.Range("H113").Value = 2370.05
$ws.Range("I113").Value = 2118.875
$ws.Range("K113").Value = 2118.875
$ws.Range("M113").Value = 51.125

$ws.Range("H136").Value = 3827.4546
$ws.Range("I136").Value = 3310.2
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 9930.599999999999
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -7380.599999999999
$ws.Range("N136").Value = -32100

$ws.Range("H137").Value = 95833.336
$ws.Range("J137").Value = 101700
$ws.Range("L137").Value = 101700
$ws.Range("N137").Value = -111900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1559.375
$ws.Range("I122").Value = 1492.4445
$ws.Range("K122").Value = 4477.333500000001
$ws.Range("M122").Value = -2027.333500000001

$ws.Range("H132").Value = 3195.4
$ws.Range("I132").Value = 3217.111
$ws.Range("K132").Value = 9651.332999999999
$ws.Range("M132").Value = -7121.332999999999
